# Rename "Device" sheet/terminology to "Apparatus" throughout the workbook
# (commit: Change "device" to "apparatus" (#13)).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Device")

# Update the cell text that referred to "Device" -> "Apparatus"
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"

# Rename the worksheet itself
$ws.Name = "Apparatus"

# Make this sheet the active / selected tab (it becomes the active tab in the
# saved workbook, matching the new bookViews/sheetViews state)
$ws.Activate()
